# Updated cryptos list on Tue Jun  6 19:57:00 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto-currency rows in the sheet. Column D holds plain text (not real
# numbers) in this sheet, so whenever the new price string would otherwise
# be re-interpreted by Excel as a genuine Number (which would silently
# drop meaningful trailing zeros, e.g. "0.9960" -> 0.996, or normalise
# "0.000008211" to scientific notation) the assignment is prefixed with a
# leading apostrophe to force literal text, matching how Excel itself
# would store a manually-typed value that looks numeric.
#
# Rows 36 & 37 also swapped rank order this run: HuobiToken now ranks above
# MXToken, so their Coin name / Link / Price / Volume values swap rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 36 & 37 swapped rank: HuobiToken now ranks above MXToken ---
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.885"
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'3.313"
$ws.Range("E37").Value = "  +9.48%  "

# --- Price / Volume(1h) refresh for all other rows ---
$ws.Range("D2").Value = "27.063.36"
$ws.Range("E2").Value = "  +5.71%  "
$ws.Range("D3").Value = "1.888.15"
$ws.Range("E3").Value = "  +4.50%  "
$ws.Range("D4").Value = "'0.9960"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'284.49"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").Value = "'0.9960"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "'0.5255"
$ws.Range("E7").Value = "  +4.38%  "
$ws.Range("D8").Value = "'0.3553"
$ws.Range("D9").Value = "'45.40"
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").Value = "'0.07103"
$ws.Range("E10").Value = "  +6.52%  "
$ws.Range("D11").Value = "'20.57"
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("D12").Value = "'0.8296"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "'0.07772"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "1.875.19"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "'5.232"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "'90.62"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "'0.9962"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'14.58"
$ws.Range("E18").Value = "  +4.91%  "
$ws.Range("D19").Value = "'0.000008211"
$ws.Range("E19").Value = "  +3.51%  "
$ws.Range("D20").Value = "'0.9980"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "27.061.34"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("D22").Value = "'4.810"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").Value = "'10.24"
$ws.Range("D24").Value = "'6.281"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").Value = "'2.436"
$ws.Range("E25").Value = "  +15.05%  "
$ws.Range("D26").Value = "'146.11"
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").Value = "'17.55"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("D28").Value = "'1.677"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "'112.54"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("D30").Value = "'4.453"
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").Value = "'4.420"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").Value = "'0.08888"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'0.04958"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("D34").Value = "'1.192"
$ws.Range("E34").Value = "  +5.89%  "
$ws.Range("D35").Value = "'0.7583"
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("D38").Value = "'2.440"
$ws.Range("E38").Value = "  +6.17%  "
$ws.Range("D39").Value = "'0.5372"
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("D40").Value = "'0.01895"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "'0.9873"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("D42").Value = "'117.36"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("D43").Value = "'6.350"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("D44").Value = "'8.315"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "'0.4682"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").Value = "'0.9960"
$ws.Range("D47").Value = "'0.1379"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'9.525"
$ws.Range("D49").Value = "'37.05"
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").Value = "'1.533"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "'0.05954"
$ws.Range("E51").Value = "  +2.47%  "
